$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 97633008
$ws.Range("B2").Value = 96354
$ws.Range("D2").Value = "LC"
$ws.Range("E2").Value = 221952
$ws.Range("F2").Value = "Spindelblomster"
$ws.Range("G2").Value = "Neottia cordata"
$ws.Range("H2").Value = "(L.) Rich."
$ws.Range("Q2").Value = 398956.7786445262
$ws.Range("R2").Value = 6788881.850361901

$ws.Range("A3").Value = 97633011
$ws.Range("B3").Value = 76863
$ws.Range("D3").Value = "VU"
$ws.Range("E3").Value = 498
$ws.Range("F3").Value = "Liten sotlav"
$ws.Range("G3").Value = "Acolium karelicum"
$ws.Range("H3").Value = "(Vain.) M.Prieto & Wedin"
$ws.Range("Q3").Value = 398704.3123741738
$ws.Range("R3").Value = 6788848.641365388

$ws.Range("A4").Value = 97633015
$ws.Range("B4").Value = 89406
$ws.Range("D4").Value = "NT"
$ws.Range("E4").Value = 1204
$ws.Range("F4").Value = "Gränsticka"
$ws.Range("G4").Value = "Phellopilus nigrolimitatus"
$ws.Range("H4").Value = "(Romell) Niemelä, T.Wagner & M.Fisch."
$ws.Range("Q4").Value = 398672.5184773419
$ws.Range("R4").Value = 6788901.13211614

$ws.Range("A5").Value = 97633014
$ws.Range("B5").Value = 94121
$ws.Range("D5").Value = "NT"
$ws.Range("E5").Value = 53
$ws.Range("F5").Value = "Vedtrappmossa"
$ws.Range("G5").Value = "Crossocalyx hellerianus"
$ws.Range("H5").Value = "(Nees ex Lindenb.) Meyl."
$ws.Range("Q5").Value = 398696.0986209051
$ws.Range("R5").Value = 6788848.396463233

$ws.Range("A6").Value = 97633007
$ws.Range("B6").Value = 77506
$ws.Range("D6").Value = "NT"
$ws.Range("E6").Value = 6425
$ws.Range("F6").Value = "Garnlav"
$ws.Range("G6").Value = "Alectoria sarmentosa"
$ws.Range("H6").Value = "(Ach.) Ach."
$ws.Range("Q6").Value = 398646.143454886
$ws.Range("R6").Value = 6788857.068735377

$ws.Range("A7").Value = 97633016
$ws.Range("B7").Value = 78603
$ws.Range("D7").Value = "LC"
$ws.Range("E7").Value = 6464
$ws.Range("F7").Value = "Luddlav"
$ws.Range("G7").Value = "Nephroma resupinatum"
$ws.Range("H7").Value = "(L.) Ach."
$ws.Range("Q7").Value = 398669.6522660756
$ws.Range("R7").Value = 6788902.178847552

$ws.Range("A8").Value = 97633012
$ws.Range("B8").Value = 89406
$ws.Range("D8").Value = "NT"
$ws.Range("E8").Value = 1204
$ws.Range("F8").Value = "Gränsticka"
$ws.Range("G8").Value = "Phellopilus nigrolimitatus"
$ws.Range("H8").Value = "(Romell) Niemelä, T.Wagner & M.Fisch."
$ws.Range("Q8").Value = 398699.9434737806
$ws.Range("R8").Value = 6788847.803485386

$ws.Range("A9").Value = 97633013
$ws.Range("B9").Value = 89392
$ws.Range("D9").Value = "NT"
$ws.Range("E9").Value = 1202
$ws.Range("F9").Value = "Ullticka"
$ws.Range("G9").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H9").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q9").Value = 398510.0508583009
$ws.Range("R9").Value = 6788941.975751169
